# Update column F (dSF) values to match repulled data / push all data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    4  = -4
    5  = -4
    8  = -4
    10 = -3
    11 = -7
    12 = -3
    13 = 3
    14 = 1
    15 = -5
    16 = 3
    17 = -5
    18 = 2
    19 = -3
    20 = 2
    22 = -6
    23 = -2
    24 = 11
    25 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
